# Updates the cryptos list values (prices / 1h volume %) and fixes a few
# row re-orderings, as scraped on Mon Sep 25 15:30:30 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping the cell as plain
# text (matches the source feed, which always stores Price as text),
# and make sure no stray number-format style sticks around afterwards.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple value updates (Price / Volume(1h) columns) ---

$ws.Range("D2").Value = "26.235.53"
$ws.Range("E2").Value = "  -1.86%  "

$ws.Range("D3").Value = "1.583.39"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  -0.64%  "

Set-TextValue $ws.Range("D5") "209.11"
$ws.Range("E5").Value = "  -1.31%  "

Set-TextValue $ws.Range("D6") "0.500"
$ws.Range("E6").Value = "  -2.61%  "

Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.73%  "

Set-TextValue $ws.Range("D8") "0.0611"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("E9").Value = "  -0.78%  "

Set-TextValue $ws.Range("D10") "19.56"
$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").Value = "1.801.89"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "1.561.37"
$ws.Range("E14").Value = "  -2.51%  "

Set-TextValue $ws.Range("D15") "0.517"
$ws.Range("E15").Value = "  -1.53%  "

Set-TextValue $ws.Range("D16") "64.41"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").Value = "26.231.33"
$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -1.60%  "

$ws.Range("E19").Value = "  +2.54%  "

Set-TextValue $ws.Range("D20") "211.20"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("E22").Value = "  -1.05%  "

Set-TextValue $ws.Range("D23") "2.17"
$ws.Range("E23").Value = "  -2.05%  "

Set-TextValue $ws.Range("D24") "8.84"
$ws.Range("E24").Value = "  -2.10%  "

Set-TextValue $ws.Range("D25") "144.23"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("E28").Value = "  -1.64%  "

Set-TextValue $ws.Range("D29") "15.29"
$ws.Range("E29").Value = "  -0.33%  "

Set-TextValue $ws.Range("D30") "0.0507"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("E32").Value = "  -1.37%  "

Set-TextValue $ws.Range("D33") "2.99"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").Value = "1.284.75"
$ws.Range("E34").Value = "  +0.33%  "

Set-TextValue $ws.Range("D35") "2.44"
$ws.Range("E35").Value = "  -1.73%  "

Set-TextValue $ws.Range("D36") "0.605"
$ws.Range("E36").Value = "  +2.50%  "

$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("E38").Value = "  -8.44%  "

Set-TextValue $ws.Range("D39") "0.0167"
$ws.Range("E39").Value = "  -1.56%  "

Set-TextValue $ws.Range("D40") "0.812"
$ws.Range("E40").Value = "  -1.57%  "

Set-TextValue $ws.Range("D41") "0.999"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("E42").Value = "  +2.28%  "

# --- Rows 43-45 get reshuffled (TrustWalletToken / MXToken / Aave) ---

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D43") "2.13"
$ws.Range("E43").Value = "  -2.73%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D44") "62.53"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D45") "0.762"
$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("D46").Value = "1.717.58"
$ws.Range("E46").Value = "  -1.24%  "

Set-TextValue $ws.Range("D47") "88.66"
$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("E48").Value = "  -3.67%  "

# --- Rows 49-50 get reshuffled (BabyDogeCoin / Algorand) ---

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.100"
$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0988"
$ws.Range("E50").Value = "  -6.72%  "

$ws.Range("E51").Value = "  -1.43%  "
